$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INS")

# -------------------------------------------------------------------
# 1) Growth-rate correction for the existing "T-CAR-ICE_B*" row (F13/F14)
#    1.15 -> 1.17
# -------------------------------------------------------------------
$ws.Cells.Item(14,6).Formula = "=F13*(1.17^32)"

# -------------------------------------------------------------------
# 2) New capacity-bound block for "T-MGT-BEV*" (two 3-row groups,
#    rows 25-27 and 28-30), mirroring the existing T-LGT-BEV*/T-LGT-FCV*
#    blocks in rows 19-21 / 22-24.
# -------------------------------------------------------------------

function Copy-RowFormat($srcRow, $dstRow, $firstCol, $lastCol) {
    $src = $ws.Range($ws.Cells.Item($srcRow, $firstCol), $ws.Cells.Item($srcRow, $lastCol))
    $dst = $ws.Range($ws.Cells.Item($dstRow, $firstCol), $ws.Cells.Item($dstRow, $lastCol))
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# --- Row 25 (pattern of row 19): C,D,E,F,H,I ---
Copy-RowFormat 19 25 3 4
Copy-RowFormat 19 25 8 9
$ws.Cells.Item(25,3).Value = "UP"
$ws.Cells.Item(25,4).Value = "CAP_BND"
$ws.Cells.Item(25,5).Value = 2019
$ws.Cells.Item(25,6).Value = 0.4
$ws.Cells.Item(25,8).Value = "T-MGT-BEV*"
$ws.Cells.Item(25,9).Value = "*New"

# --- Row 26 (pattern of row 20): C,D,E,F,H,I ---
Copy-RowFormat 20 26 3 4
Copy-RowFormat 20 26 6 6
Copy-RowFormat 20 26 8 9
$ws.Cells.Item(26,3).Value = "UP"
$ws.Cells.Item(26,4).Value = "CAP_BND"
$ws.Cells.Item(26,5).Value = 2050
$ws.Cells.Item(26,6).Formula = "=F25*(1.15^32)"
$ws.Cells.Item(26,8).Value = "T-MGT-BEV*"
$ws.Cells.Item(26,9).Value = "*New"

# --- Row 27 (pattern of row 21): B..N ---
Copy-RowFormat 21 27 2 14
$ws.Cells.Item(27,3).Value = "UP"
$ws.Cells.Item(27,4).Value = "CAP_BND"
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 5
$ws.Cells.Item(27,8).Value = "T-MGT-BEV*"
$ws.Cells.Item(27,9).Value = "*New"

# --- Row 28 (pattern of row 22): C,D,E,F,H,I ---
Copy-RowFormat 22 28 3 4
Copy-RowFormat 22 28 8 9
$ws.Cells.Item(28,3).Value = "UP"
$ws.Cells.Item(28,4).Value = "CAP_BND"
$ws.Cells.Item(28,5).Value = 2019
$ws.Cells.Item(28,6).Value = 0.4
$ws.Cells.Item(28,8).Value = "T-MGT-BEV*"
$ws.Cells.Item(28,9).Value = "*New"

# --- Row 29 (pattern of row 23): C,D,E,F,H,I ---
Copy-RowFormat 23 29 3 4
Copy-RowFormat 23 29 6 6
Copy-RowFormat 23 29 8 9
$ws.Cells.Item(29,3).Value = "UP"
$ws.Cells.Item(29,4).Value = "CAP_BND"
$ws.Cells.Item(29,5).Value = 2050
$ws.Cells.Item(29,6).Formula = "=F28*(1.15^32)"
$ws.Cells.Item(29,8).Value = "T-MGT-BEV*"
$ws.Cells.Item(29,9).Value = "*New"

# --- Row 30 (pattern of row 24): B..N ---
Copy-RowFormat 24 30 2 14
$ws.Cells.Item(30,3).Value = "UP"
$ws.Cells.Item(30,4).Value = "CAP_BND"
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 5
$ws.Cells.Item(30,8).Value = "T-MGT-BEV*"
$ws.Cells.Item(30,9).Value = "*New"

# -------------------------------------------------------------------
# 3) The "T-LGT-BEV*"/"T-LGT-FCV*" capacity-bound rows (21/24) lose
#    their now-orphaned direct-format ("Normal 14" font) and fall back
#    to the same plain border style already used by the rest of the
#    block (mirrors the sibling cells I21/I24 etc.).
# -------------------------------------------------------------------
Copy-RowFormat 9 21 8 8
Copy-RowFormat 9 24 8 8
$ws.Cells.Item(21,8).Value = "T-LGT-BEV*"
$ws.Cells.Item(24,8).Value = "T-LGT-FCV*"

# -------------------------------------------------------------------
# 4) Update the active selection to mirror the post-edit state
# -------------------------------------------------------------------
$ws.Range("F29").Select()

Write-Host "done"
